# config.xlsx update:
#  - Add a new "Test Suite" row (value "Smoke") above the existing "Browser"
#    row on both the "UAT Apollo" and "UAT ATeam" sheets.
#  - The relocated "Browser" row's value becomes "Firefox" (the new default)
#    on sheets where it wasn't already.
#  - Data validation list for the browser cell moves down with it; a new
#    validation (Smoke/Function) is added for the Test Suite cell.
#  - Hyperlink on the "Login Main Site" row moves down one row as well.

$wb = $excel.ActiveWorkbook

$xlPasteValues  = -4163
$xlPasteFormats = -4122

function Shift-ConfigSheetDown {
    param($ws)

    # Push rows 9..2 down to 10..3 (bottom-up so we never clobber data we
    # still need to read), carrying both the literal value and the cell
    # formatting (border/number format/etc.) along with it.
    for ($r = 9; $r -ge 2; $r--) {
        $src = $ws.Range("A" + $r + ":C" + $r)
        $dst = $ws.Range("A" + ($r + 1) + ":C" + ($r + 1))

        $src.Copy()
        $dst.PasteSpecial($xlPasteValues)
        $src.Copy()
        $dst.PasteSpecial($xlPasteFormats)
    }
    $excel.CutCopyMode = $false

    # New header-matching row 2: "Test Suite" / "Smoke", styled the same as
    # the Env sheet's data row (thin border, no bold/left-align numfmt).
    $ws.Range("A2").Value = "Test Suite"
    $ws.Range("B2").Value = "Smoke"

    $envWs = $wb.Worksheets.Item("Env")
    $envWs.Range("A2:C2").Copy()
    $ws.Range("A2:C2").PasteSpecial($xlPasteFormats)
    $excel.CutCopyMode = $false

    # The "Browser" row (now row 3) should default to Firefox.
    $ws.Range("B3").Value = "Firefox"

    # Move the hyperlink that was on the login-url row (row 7) down to its
    # new home at row 8.
    $linkAddress = $ws.Range("B7").Hyperlinks.Item(1).Address
    $ws.Range("B7").Hyperlinks.Delete()
    $ws.Hyperlinks.Add($ws.Range("B8"), $linkAddress)

    # Data validation: the old Firefox/Chrome list for the Browser cell
    # moves from B2 to B3; a fresh Smoke/Function list is added for the new
    # Test Suite cell at B2.
    $ws.Range("B2").Validation.Delete()
    $ws.Range("B3").Validation.Add(3, 1, 1, """Firefox, Chrome""")
    $ws.Range("B2").Validation.Add(3, 1, 1, """Smoke, Function""")
}

$wsApollo = $wb.Worksheets.Item("UAT Apollo")
$wsATeam  = $wb.Worksheets.Item("UAT ATeam")

Shift-ConfigSheetDown $wsApollo
Shift-ConfigSheetDown $wsATeam

# Restore the selections Excel would have left behind after the edits.
$wsEnv = $wb.Worksheets.Item("Env")
$wsEnv.Select()
$wsEnv.Rows.Item(3).Select()

$wsATeam.Select()
$wsATeam.Range("A2:C2").Select()

$wsApollo.Select()
$wsApollo.Range("B2").Select()
